# Regen sval data to filter save games
# Update the numeric stat columns (B:E) and the derived sum column (G)
# for rows 2-5 on the active sheet. Column F (Win) stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.3464964993005633; C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.896700893398075 }
    3 = @{ B = 0.1554434735375247; C = 1.65323645889881;   D = 0.7127328510149897; E = 6.48142807727062;   G = 9.002840860721944 }
    4 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 6.741336633845642 }
    5 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 6.048734245549538 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
